$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.115.01"
$ws.Range("D3").Value = "2.421.83"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'553.85"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").Value = "'138.68"
$ws.Range("E6").Value = "  +3.09%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("D10").Value = "'5.78"
$ws.Range("D11").Value = "'0.359"
$ws.Range("E11").Value = "  +0.69%  "
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").Value = "'24.95"
$ws.Range("E13").Value = "  +4.92%  "
$ws.Range("D14").Value = "2.853.71"
$ws.Range("E14").Value = "  +3.12%  "
$ws.Range("D15").Value = "60.026.56"
$ws.Range("E15").Value = "  +3.84%  "
$ws.Range("E16").Value = "  +3.33%  "
$ws.Range("D17").Value = "2.422.74"
$ws.Range("E17").Value = "  +3.38%  "
$ws.Range("E18").Value = "  +6.41%  "
$ws.Range("D19").Value = "'4.39"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Value = "'332.21"
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D23").Value = "'65.14"
$ws.Range("E23").Value = "  +3.96%  "
$ws.Range("E24").Value = "  +3.58%  "
$ws.Range("D25").Value = "'8.59"
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +0.46%  "
$ws.Range("D28").Value = "0.0₃0790"
$ws.Range("E28").Value = "  +7.44%  "
$ws.Range("E29").Value = "  +1.62%  "
$ws.Range("D30").Value = "'6.29"
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("D31").Value = "'169.82"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("D33").Value = "'18.71"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D35").Value = "'1.30"
$ws.Range("E35").Value = "  +5.60%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").Value = "'4.22"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  +11.60%  "
$ws.Range("D40").Value = "'320.93"
$ws.Range("E40").Value = "  +11.30%  "
$ws.Range("D41").Value = "'39.43"
$ws.Range("E41").Value = "  +0.97%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("D43").Value = "'140.11"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'0.0961"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "'0.0521"
$ws.Range("E45").Value = "  +2.55%  "
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("D47").Value = "'0.413"
$ws.Range("E47").Value = "  +8.28%  "
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "'17.78"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("E51").Value = "  -0.23%  "
